# Restore C10 ("R30" rule's "Integer min" / C1 value) from 18 to 1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
